$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($r = 2; $r -le 515; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
